# TC-63769 / NGC-1826 test data update
# - Sheet "Add Devices Loop A": update Description/UsedForMethods/UserStory cells
#   and the DC-units numbers in G1/G2; fix the selected range.
# - Sheet "Update Devices": populate the same Description/UsedForMethods/UserStory
#   cells (previously blank) and update the DC-units numbers in G1/G2; fix the
#   selected range and drop the stale topLeftCell scroll position.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Add Devices Loop A"
$ws2 = $wb.Worksheets.Item(2)   # "Update Devices"

# --- Sheet 1: "Add Devices Loop A" ---------------------------------------
$ws1.Range("B4").Value = "NGC-1826/TC-63769"
$ws1.Range("B3").Value = "VerifyCurrentDCCalculation"
$ws1.Range("B2").Value = "Verify current DCCalculation "

$ws1.Range("G1").Value = 345.4
$ws1.Range("G2").Value = 332.3

$ws1.Rows.Item(2).AutoFit()

$ws1.Range("B2:B4").Select()

# --- Sheet 2: "Update Devices" -------------------------------------------
$ws2.Range("B4").Value = "NGC-1826/TC-63769"
$ws2.Range("B3").Value = "VerifyCurrentDCCalculation"
$ws2.Range("B2").Value = "Verify current DCCalculation "

$ws2.Range("G1").Value = 345.4
$ws2.Range("G2").Value = 332.3

$ws2.Range("B2:B4").Select()
